$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'25.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("B4").Value = "LEO"
$ws.Range("C4").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D4").Value = "'3.500"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "3LEOLEO"
$ws.Range("B5").Value = "HuobiToken"
$ws.Range("C5").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D5").Value = "'5.178"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "4HuobiTokenHT"
$ws.Range("B6").Value = "Cronos"
$ws.Range("C6").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D6").Value = "'0.05708"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "5CronosCRO"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'6.481"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "6KuCoinTokenKCS"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'2.968"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "7GateTokenGT"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.8101"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "8MXTokenMX"
$ws.Range("B10").Value = "FTXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D10").Value = "'0.8338"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9FTXTokenFTT"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1331"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.06960"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.02830"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09377"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").Value = "'0.0005999"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "15OneONEWorstin24h"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006114"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "16TigerCashTCH"
$ws.Range("D18").Value = "'2.123"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Value = "'0.03202"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.1338"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'3.740"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.04684"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Value = "'0.001235"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.004243"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.00009698"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001744"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "27UpBotsUBXT"
$ws.Range("D40").Value = "'0.03630"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006289"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1048"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.003000"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.007359"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005282"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.2000"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.002295"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("D50").Style = "Normal"
